$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 (сасмадрас): add party names for левмил, консерватизм, соцдем, эгалитаризм, эгоизм
$ws.Range("J7").Value = "Коммунистическая партия"
$ws.Range("K7").Value = "Единый Пайноплус"
$ws.Range("L7").Value = "Соц.-дем. Партия"
$ws.Range("M7").Value = "Совет крестьянских общин"
$ws.Range("N7").Value = "Союз эгоистов"

# Row 24 (арванта) / Row 27 (велбурч): add party names in column G (социализм)
# Shared-string order requires G27 to be entered before G24
$ws.Range("G27").Value = "Левая Интеелегенция"
$ws.Range("G24").Value = "Партия Равенства"

# Update selection / view to match the authored state
$ws.Range("G23").Select()
